$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking Price values in column D are written as
# plain text (matching the workbook's existing inline-string convention)
# instead of being auto-converted to numbers by Excel.
$priceCells = @("D2","D3","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D23","D24","D27","D40","D41","D42","D43","D44","D47")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Price (column D) updates ---
$ws.Range("D2").Value = "246.11"
$ws.Range("D3").Value = "22.06"
$ws.Range("D5").Value = "0.05778"
$ws.Range("D7").Value = "6.322"
$ws.Range("D8").Value = "0.8176"
$ws.Range("D9").Value = "0.9473"
$ws.Range("D10").Value = "0.01121"
$ws.Range("D11").Value = "0.1428"
$ws.Range("D12").Value = "0.07508"
$ws.Range("D13").Value = "0.03140"
$ws.Range("D14").Value = "0.03000"
$ws.Range("D15").Value = "4.157"
$ws.Range("D16").Value = "0.09407"
$ws.Range("D17").Value = "0.001595"
$ws.Range("D18").Value = "0.04816"
$ws.Range("D19").Value = "0.006200"
$ws.Range("D20").Value = "0.004124"
$ws.Range("D21").Value = "0.0009969"
$ws.Range("D23").Value = "3.776"
$ws.Range("D24").Value = "2.225"
$ws.Range("D27").Value = "0.0004000"
$ws.Range("D40").Value = "0.03897"
$ws.Range("D41").Value = "0.006371"
$ws.Range("D42").Value = "0.1075"
$ws.Range("D43").Value = "0.003001"
$ws.Range("D44").Value = "0.006551"
$ws.Range("D47").Value = "0.3801"

# --- Row 9: FTXToken volume label correction ---
$ws.Range("E9").Value = "8FTXTokenFTT"

# --- Rows 10-18: new "One" entry inserted, remaining rows shifted down ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E16").Value = "15BitMartTokenBMX"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E18").Value = "17CoinExTokenCET"
